$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 378; existing rows 378..447 shift down to 379..448.
$ws.Rows.Item(378).Insert()

# Populate the newly inserted row 378 with the new record
# (same Mercado/Categoria data as its neighbours, new Fecha/Volumen/Precio values).
$ws.Range("A378").Value = 10
$ws.Range("B378").Value = "Vega Modelo de Temuco"
$ws.Range("C378").Value = "La Araucanía"
$ws.Range("D378").Value = 45258
$ws.Range("E378").Value = 9
$ws.Range("F378").Value = 100112039
$ws.Range("G378").Value = "Ciboulette"
$ws.Range("H378").Value = "Sin especificar"
$ws.Range("I378").Value = "Primera"
$ws.Range("J378").Value = 35
$ws.Range("K378").Value = 7000
$ws.Range("L378").Value = 7000
$ws.Range("M378").Value = 7000
$ws.Range("N378").Value = '$/docena de atados'
$ws.Range("O378").Value = "Provincia de Cautín"
$ws.Range("P378").Value = 2333
$ws.Range("Q378").Value = 3
$ws.Range("R378").Value = "Hortaliza"
